# Apply "contingencies with rene fine": add two new line rows (line7, line8)
# after line6, shifting the existing extr1..extr8 rows down by two rows, and
# update a handful of numeric/boolean values along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final state for data rows 2..17 (columns A-E): index, name, from_bus, to_bus, in_service
$data = @(
  @(0,  "line1", 7,  9,  $false),
  @(1,  "line2", 9,  8,  $true),
  @(2,  "line3", 8,  10, $true),
  @(3,  "line4", 8,  11, $true),
  @(4,  "line5", 10, 5,  $true),
  @(5,  "line6", 12, 8,  $true),
  @(6,  "line7", 14, 11, $true),
  @(7,  "line8", 16, 9,  $true),
  @(8,  "extr1", 5,  12, $true),
  @(9,  "extr2", 5,  9,  $true),
  @(10, "extr3", 10, 11, $false),
  @(11, "extr4", 7,  8,  $true),
  @(12, "extr5", 9,  11, $true),
  @(13, "extr6", 7,  11, $false),
  @(14, "extr7", 5,  7,  $false),
  @(15, "extr8", 8,  5,  $true)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row = $row + 1
}

# Rows 16 and 17 are brand new (previously the table ended at row 15), so
# column A needs the same bold/centered/bordered style used by the rest
# of the "name index" column. Copy that formatting down from row 15.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
